$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The journal's date column used "MM/DD/YY;@" - switch it to "DD/MM/YY;@"
$ws.Range("A5:A31").NumberFormat = "DD/MM/YY;@"

# Fill in the new journal entry on row 13
$ws.Range("A13").Value = 43178
$ws.Range("B13").Value = "Réalisation de la main frame"
$ws.Range("C13").Value = 6
$ws.Rows.Item(13).RowHeight = 13.8

# Move the active selection to A14
$ws.Range("A14").Select()
